# Generate Report for Archive
#
# The localization run finished and the sheets that used to say
# "Ready for handoff" should now read "In Translation" while the
# report is archived. The two narrow "status" columns that show this
# value on each sheet are also re-sized to better fit the shorter text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
$newWidth  = 13.4101845877511

# --- Overview sheet: zh-cn / de-de status columns are E and F ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F4").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

# --- zh-cn sheet: Status column is C ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C4").Value = $newStatus
$zhcn.Columns.Item(3).ColumnWidth = $newWidth

# --- de-de sheet: Status column is C ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C4").Value = $newStatus
$dede.Columns.Item(3).ColumnWidth = $newWidth
